$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the four DDL statements in column H (rows 29-32) to include " engine=BTREE"
$ws.Range("H29").Value = "create database MYBTREE_DDL_028;create table MYBTREE_DDL_028.MYBTREEDDL028_TBL01(id int, name varchar(20), primary key(id)) engine=BTREE"
$ws.Range("H30").Value = "create database MYBTREE_DDL_029;create table MYBTREE_DDL_029.MYBTREEDDL029_TBL01(id int, name varchar(20), primary key(id)) engine=BTREE;drop table MYBTREE_DDL_029.MYBTREEDDL029_TBL01"
$ws.Range("H31").Value = "create table MYSQL.MYBTREEDDL030_TBL01(id int not null auto_increment, name varchar(20), primary key(id)) engine=BTREE"
$ws.Range("H32").Value = "create table MYSQL.MYBTREEDDL031_TBL01(id int not null auto_increment, name varchar(20), primary key(id)) engine=BTREE;drop table MYSQL.MYBTREEDDL031_TBL01"

# Scroll the view and move the active selection to reflect the saved view state
$ws.Activate()
$ws.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G32").Select()
